# The "Förändrad" (changed) date column (C) was bumped from 2023-09-21
# (serial 45190) to 2023-09-23 (serial 45192) for every data row
# (rows 2-267) on the single worksheet "Avverkningsanmälningar".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C267").Value = 45192
